# The workbook's active/selected tab is "Sheet7" (physically xl/worksheets/sheet1.xml)
# and the second tab is "Sheet5" (physically xl/worksheets/sheet2.xml).
#
# The edit removes the "KA BIAS" data row (row 27) from the "Sheet7" sheet,
# which shifts every following row up by one (dimension A1:E38 -> A1:E37) and
# drops one shared-string reference (sharedStrings count 287 -> 286). It also
# leaves the sheets with updated cursor/selection positions: "A10" on Sheet7
# and "A35" on Sheet5.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet7")
$ws2 = $wb.Worksheets.Item("Sheet5")

# Delete the whole row 27 (the "KA BIAS" record) on the active sheet; this
# shifts rows 28:38 up to 27:37 and recalculates the used range/dimension.
$ws1.Rows("27:27").Delete()

# Update Sheet5's selection (visited while editing) without leaving it as the
# active tab.
$ws2.Activate()
$ws2.Range("A35").Select()

# Re-activate Sheet7 (keeps it the tab shown/selected on reopen) and leave the
# cursor on A10.
$ws1.Activate()
$ws1.Range("A10").Select()
